$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 25
$ws.Range("F2").Value = 31
$ws.Range("G2").Value = "adam"
$ws.Range("I2").Value = 32
$ws.Range("J2").Value = 25.32245183180559
$ws.Range("K2").Value = 1047.963919416445
$ws.Range("L2").Value = 32.37227084120676
$ws.Range("M2").Value = 0.1535140507568457
